# Apply updated cryptocurrency price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.NumberFormat = "@"
$r.Value = '25.998.44'
$r.ClearFormats()
$ws.Range('E2').Value = '  +0.36%  '
$r = $ws.Range('D3')
$r.NumberFormat = "@"
$r.Value = '1.638.66'
$r.ClearFormats()
$ws.Range('E4').Value = '  -0.32%  '
$r = $ws.Range('D5')
$r.NumberFormat = "@"
$r.Value = '215.00'
$r.ClearFormats()
$ws.Range('E5').Value = '  +0.07%  '
$r = $ws.Range('D6')
$r.NumberFormat = "@"
$r.Value = '0.5125'
$r.ClearFormats()
$ws.Range('E6').Value = '  +1.60%  '
$ws.Range('E7').Value = '  -0.19%  '
$r = $ws.Range('D8')
$r.NumberFormat = "@"
$r.Value = '0.2583'
$r.ClearFormats()
$ws.Range('E8').Value = '  +0.37%  '
$r = $ws.Range('D9')
$r.NumberFormat = "@"
$r.Value = '0.06368'
$r.ClearFormats()
$ws.Range('E9').Value = '  -0.54%  '
$r = $ws.Range('D10')
$r.NumberFormat = "@"
$r.Value = '19.81'
$r.ClearFormats()
$r = $ws.Range('D11')
$r.NumberFormat = "@"
$r.Value = '0.07796'
$r.ClearFormats()
$ws.Range('E11').Value = '  +0.12%  '
$r = $ws.Range('D12')
$r.NumberFormat = "@"
$r.Value = '4.290'
$r.ClearFormats()
$ws.Range('E12').Value = '  +0.02%  '
$r = $ws.Range('D13')
$r.NumberFormat = "@"
$r.Value = '1.645.30'
$r.ClearFormats()
$ws.Range('E13').Value = '  +0.04%  '
$r = $ws.Range('D14')
$r.NumberFormat = "@"
$r.Value = '0.5466'
$r.ClearFormats()
$ws.Range('E14').Value = '  +0.47%  '
$r = $ws.Range('D15')
$r.NumberFormat = "@"
$r.Value = '64.62'
$r.ClearFormats()
$ws.Range('E15').Value = '  -0.69%  '
$r = $ws.Range('D16')
$r.NumberFormat = "@"
$r.Value = '0.0₅7741'
$r.ClearFormats()
$ws.Range('E16').Value = '  -1.70%  '
$r = $ws.Range('D17')
$r.NumberFormat = "@"
$r.Value = '26.008.95'
$r.ClearFormats()
$ws.Range('E17').Value = '  +0.16%  '
$ws.Range('E18').Value = '  -0.23%  '
$r = $ws.Range('D19')
$r.NumberFormat = "@"
$r.Value = '197.74'
$r.ClearFormats()
$ws.Range('E19').Value = '  -0.03%  '
$ws.Range('E20').Value = '  +0.91%  '
$r = $ws.Range('D21')
$r.NumberFormat = "@"
$r.Value = '9.973'
$r.ClearFormats()
$ws.Range('E21').Value = '  +0.02%  '
$r = $ws.Range('D22')
$r.NumberFormat = "@"
$r.Value = '6.089'
$r.ClearFormats()
$ws.Range('E22').Value = '  +1.04%  '
$r = $ws.Range('D23')
$r.NumberFormat = "@"
$r.Value = '1.005'
$r.ClearFormats()
$ws.Range('E23').Value = '  -0.21%  '
$r = $ws.Range('D24')
$r.NumberFormat = "@"
$r.Value = '1.892'
$r.ClearFormats()
$ws.Range('E24').Value = '  +1.10%  '
$r = $ws.Range('D25')
$r.NumberFormat = "@"
$r.Value = '141.97'
$r.ClearFormats()
$ws.Range('E25').Value = '  +1.19%  '
$ws.Range('E26').Value = '  +7.53%  '
$ws.Range('E27').Value = '  -0.13%  '
$r = $ws.Range('D28')
$r.NumberFormat = "@"
$r.Value = '15.68'
$r.ClearFormats()
$ws.Range('E28').Value = '  -0.24%  '
$r = $ws.Range('D29')
$r.NumberFormat = "@"
$r.Value = '1.238'
$r.ClearFormats()
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('E30').Value = '  -2.37%  '
$ws.Range('E31').Value = '  +0.60%  '
$r = $ws.Range('D32')
$r.NumberFormat = "@"
$r.Value = '3.213'
$r.ClearFormats()
$ws.Range('E32').Value = '  +0.59%  '
$r = $ws.Range('D33')
$r.NumberFormat = "@"
$r.Value = '1.541'
$r.ClearFormats()
$ws.Range('E33').Value = '  +0.32%  '
$r = $ws.Range('D34')
$r.NumberFormat = "@"
$r.Value = '2.375'
$r.ClearFormats()
$ws.Range('E34').Value = '  +0.35%  '
$r = $ws.Range('D35')
$r.NumberFormat = "@"
$r.Value = '0.9154'
$r.ClearFormats()
$ws.Range('E35').Value = '  +2.39%  '
$r = $ws.Range('D36')
$r.NumberFormat = "@"
$r.Value = '2.587'
$r.ClearFormats()
$ws.Range('E36').Value = '  -0.31%  '
$r = $ws.Range('D37')
$r.NumberFormat = "@"
$r.Value = '0.5545'
$r.ClearFormats()
$ws.Range('E37').Value = '  +0.30%  '
$r = $ws.Range('D38')
$r.NumberFormat = "@"
$r.Value = '1.111.45'
$r.ClearFormats()
$ws.Range('E38').Value = '  -2.07%  '
$ws.Range('E39').Value = '  +0.85%  '
$r = $ws.Range('D40')
$r.NumberFormat = "@"
$r.Value = '1.003'
$r.ClearFormats()
$ws.Range('E40').Value = '  -0.36%  '
$r = $ws.Range('D41')
$r.NumberFormat = "@"
$r.Value = '2.533'
$r.ClearFormats()
$ws.Range('E41').Value = '  -1.13%  '
$r = $ws.Range('D42')
$r.NumberFormat = "@"
$r.Value = '5.533'
$r.ClearFormats()
$ws.Range('E42').Value = '  -2.95%  '
$r = $ws.Range('D43')
$r.NumberFormat = "@"
$r.Value = '0.8091'
$r.ClearFormats()
$ws.Range('E43').Value = '  -0.76%  '
$r = $ws.Range('D44')
$r.NumberFormat = "@"
$r.Value = '99.41'
$r.ClearFormats()
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('E45').Value = '  +0.29%  '
$r = $ws.Range('D46')
$r.NumberFormat = "@"
$r.Value = '1.776.94'
$r.ClearFormats()
$ws.Range('E46').Value = '  -0.01%  '
$r = $ws.Range('D47')
$r.NumberFormat = "@"
$r.Value = '0.4536'
$r.ClearFormats()
$ws.Range('E47').Value = '  +0.02%  '
$r = $ws.Range('D48')
$r.NumberFormat = "@"
$r.Value = '1.007'
$r.ClearFormats()
$ws.Range('E48').Value = '  +0.20%  '
$r = $ws.Range('D49')
$r.NumberFormat = "@"
$r.Value = '55.06'
$r.ClearFormats()
$ws.Range('E49').Value = '  -0.37%  '
$ws.Range('E50').Value = '  +4.05%  '
$ws.Range('E51').Value = '  +0.03%  '
